# "Added Vzw keyword driven demo"
#
# - Adds a new worksheet (named "null", matching the literal sheet name
#   produced upstream) at the end of the workbook with a small header row
#   re-using existing shared strings ("Test Parameters" / "Avner 1.14").
# - On the "iPhone-6 Avner" and "iPhone-6 Raj" sheets, attaches a new
#   hyperlink to the PASS cell in C5 (matching the pre-existing hyperlink
#   pattern already used for the other PASS/FAIL cells on those sheets),
#   while keeping the cell's original green "PASS" fill/format instead of
#   Excel's default blue/underlined hyperlink styling.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "null" worksheet at the end of the workbook
# ---------------------------------------------------------------------
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "null"

$newSheet.Range("A1").Value = "Test Parameters"
$newSheet.Range("B1").Value = "Avner 1.14"

# ---------------------------------------------------------------------
# 2. iPhone-6 Avner (sheet2.xml) - hyperlink on C5, keep PASS formatting
# ---------------------------------------------------------------------
$wsAvner = $wb.Worksheets.Item("iPhone-6 Avner")
$wsAvner.Hyperlinks.Add($wsAvner.Range("C5"), "C:/Users/AvnerG/git/Beton/Beton/test-output/screenshots-tests/vzw-keyword-driven-demo.png")
$wsAvner.Range("C4").Copy()
$wsAvner.Range("C5").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3. iPhone-6 Raj (sheet3.xml) - hyperlink on C5, keep PASS formatting
# ---------------------------------------------------------------------
$wsRaj = $wb.Worksheets.Item("iPhone-6 Raj")
$wsRaj.Hyperlinks.Add($wsRaj.Range("C5"), "C:/Users/AvnerG/git/Beton/Beton/test-output/screenshots-tests/vzw-keyword-driven-demo.png")
$wsRaj.Range("C4").Copy()
$wsRaj.Range("C5").PasteSpecial(-4122)
